$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tiny float-precision correction to the existing row 4 date/time value.
$ws.Range("A4").Value2 = 44317.77362215046

# New row 5 with the latest day's job-number data.
$ws.Range("A5").Value2 = 44318.77224265195
$ws.Range("B5").Value2 = 70419
$ws.Range("C5").Value2 = 59354
$ws.Range("D5").Value2 = 3170
$ws.Range("E5").Value2 = 1954
$ws.Range("F5").Value2 = 1387
$ws.Range("G5").Value2 = 18443
$ws.Range("H5").Value2 = 1347
$ws.Range("I5").Value2 = 791
$ws.Range("J5").Value2 = 196

# Column A uses a date/time display format (same format as the rest of the column).
$ws.Range("A5").NumberFormat = $ws.Range("A4").NumberFormat
